# Edit: insert two new weekly price rows for Mango (Agrícola del Norte S.A. de Arica)
# at row 132, pushing the existing data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 132, shifting existing rows 132.. down to 134..
$ws.Rows("132:133").Insert()

# --- New row 132 ---
$ws.Cells.Item(132, 1).Value = 1
$ws.Cells.Item(132, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(132, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(132, 4).Value = 44977
$ws.Cells.Item(132, 5).Value = 15
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100108
$ws.Cells.Item(132, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(132, 9).Value = 100108002
$ws.Cells.Item(132, 10).Value = "Mango"
$ws.Cells.Item(132, 11).Value = "Sin especificar"
$ws.Cells.Item(132, 12).Value = "Primera"
$ws.Cells.Item(132, 13).Value = 650
$ws.Cells.Item(132, 14).Value = 4500
$ws.Cells.Item(132, 15).Value = 5000
$ws.Cells.Item(132, 16).Value = 4808
$ws.Cells.Item(132, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(132, 18).Value = "Perú"
$ws.Cells.Item(132, 19).Value = 1202
$ws.Cells.Item(132, 20).Value = 4

# --- New row 133 ---
$ws.Cells.Item(133, 1).Value = 1
$ws.Cells.Item(133, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(133, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(133, 4).Value = 44977
$ws.Cells.Item(133, 5).Value = 15
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100108
$ws.Cells.Item(133, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(133, 9).Value = 100108002
$ws.Cells.Item(133, 10).Value = "Mango"
$ws.Cells.Item(133, 11).Value = "Sin especificar"
$ws.Cells.Item(133, 12).Value = "Segunda"
$ws.Cells.Item(133, 13).Value = 800
$ws.Cells.Item(133, 14).Value = 4500
$ws.Cells.Item(133, 15).Value = 5000
$ws.Cells.Item(133, 16).Value = 4688
$ws.Cells.Item(133, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(133, 18).Value = "Perú"
$ws.Cells.Item(133, 19).Value = 1172
$ws.Cells.Item(133, 20).Value = 4

# Make sure the date cells keep the same number format as the rest of column D
$ws.Range("D132:D133").NumberFormat = $ws.Range("D134").NumberFormat
